$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a number by Excel
# (e.g. "325.30" or "1.008") are forced to Text format first, then the
# cell style is reset to Normal afterwards so no stray formatting is left behind.
$textCells = @(
    "D4", "D5", "D7", "D8", "D9", "D10", "D11", "D13",
    "D14", "D15", "D16", "D17", "D18", "D19", "D22", "D23",
    "D24", "D26", "D27", "D28", "D29", "D30", "D31", "D32",
    "D33", "D34", "D35", "D37", "D38", "D39", "D40", "D42",
    "D43", "D44", "D46", "D47", "D48", "D49", "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values (prices in column D, 1h volume % in column E).
$ws.Range("D2").Value = "29.464.24"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "1.916.27"
$ws.Range("E3").Value = "  +1.13%  "
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.73%  "
$ws.Range("D5").Value = "325.30"
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("D7").Value = "0.4830"
$ws.Range("E7").Value = "  +1.10%  "
$ws.Range("D8").Value = "0.4078"
$ws.Range("E8").Value = "  +0.61%  "
$ws.Range("D9").Value = "0.08216"
$ws.Range("E9").Value = "  +2.51%  "
$ws.Range("D10").Value = "1.017"
$ws.Range("E10").Value = "  +1.62%  "
$ws.Range("D11").Value = "23.50"
$ws.Range("E11").Value = "  +0.81%  "
$ws.Range("D12").Value = "1.927.89"
$ws.Range("E12").Value = "  +0.52%  "
$ws.Range("D13").Value = "6.064"
$ws.Range("E13").Value = "  +2.31%  "
$ws.Range("D14").Value = "7.230"
$ws.Range("E14").Value = "  +2.58%  "
$ws.Range("D15").Value = "91.17"
$ws.Range("E15").Value = "  +1.97%  "
$ws.Range("D16").Value = "0.06807"
$ws.Range("E16").Value = "  +2.24%  "
$ws.Range("D17").Value = "1.008"
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("D18").Value = "0.00001039"
$ws.Range("E18").Value = "  +1.41%  "
$ws.Range("D19").Value = "17.72"
$ws.Range("E19").Value = "  +0.85%  "
$ws.Range("E20").Value = "  +0.67%  "
$ws.Range("D21").Value = "29.478.85"
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("D22").Value = "5.652"
$ws.Range("E22").Value = "  +2.40%  "
$ws.Range("D23").Value = "11.82"
$ws.Range("E23").Value = "  +1.21%  "
$ws.Range("D24").Value = "2.177"
$ws.Range("E24").Value = "  +1.13%  "
$ws.Range("D25").Value = "2.162.60"
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("D26").Value = "6.653"
$ws.Range("E26").Value = "  +12.08%  "
$ws.Range("D27").Value = "156.23"
$ws.Range("E27").Value = "  +1.52%  "
$ws.Range("D28").Value = "20.08"
$ws.Range("E28").Value = "  +1.78%  "
$ws.Range("D29").Value = "2.115"
$ws.Range("E29").Value = "  +1.52%  "
$ws.Range("D30").Value = "120.46"
$ws.Range("E30").Value = "  +2.22%  "
$ws.Range("D31").Value = "1.022"
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("D32").Value = "0.09573"
$ws.Range("E32").Value = "  +1.38%  "
$ws.Range("D33").Value = "5.652"
$ws.Range("E33").Value = "  +5.60%  "
$ws.Range("D34").Value = "3.549"
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("D35").Value = "1.374"
$ws.Range("E35").Value = "  -0.57%  "
$ws.Range("E36").Value = "  +1.91%  "
$ws.Range("D37").Value = "0.06115"
$ws.Range("E37").Value = "  +1.26%  "
$ws.Range("D38").Value = "1.182"
$ws.Range("E38").Value = "  +1.05%  "
$ws.Range("D39").Value = "0.5981"
$ws.Range("E39").Value = "  +2.19%  "
$ws.Range("D40").Value = "8.055"
$ws.Range("E40").Value = "  +2.60%  "
$ws.Range("E41").Value = "  +7.37%  "
$ws.Range("D42").Value = "0.1850"
$ws.Range("E42").Value = "  +0.65%  "
$ws.Range("D43").Value = "2.427"
$ws.Range("E43").Value = "  +1.20%  "
$ws.Range("D44").Value = "0.07627"
$ws.Range("E44").Value = "  -0.92%  "
$ws.Range("E45").Value = "  -3.38%  "
$ws.Range("D46").Value = "12.41"
$ws.Range("E46").Value = "  +2.22%  "
$ws.Range("D47").Value = "0.5587"
$ws.Range("E47").Value = "  +1.92%  "
$ws.Range("D48").Value = "1.958"
$ws.Range("E48").Value = "  +2.19%  "
$ws.Range("D49").Value = "117.74"
$ws.Range("E49").Value = "  +4.10%  "
$ws.Range("E50").Value = "  +4.21%  "
$ws.Range("D51").Value = "72.40"
$ws.Range("E51").Value = "  +1.18%  "

# Restore plain/default styling on the cells we forced to Text format so their
# appearance matches the rest of the sheet (no leftover custom number format).
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
